$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" property to the new publication date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: add a new mapping column (AL) for the business mapping ---
$elements = $wb.Worksheets.Item("Elements")

# New header in AL1, matching the "Mapping: RIM Mapping" header style in AK1
$elements.Range("AK1").Copy()
$elements.Range("AL1").PasteSpecial(-4122)
$elements.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR LocationStatus"

# Data rows: copy formatting from the neighboring AK column, most rows stay blank,
# only the Extension.value[x] row (row 6) carries the new mapping text "statut"
$elements.Range("AK2").Copy()
$elements.Range("AL2").PasteSpecial(-4122)

$elements.Range("AK3").Copy()
$elements.Range("AL3").PasteSpecial(-4122)

$elements.Range("AK4").Copy()
$elements.Range("AL4").PasteSpecial(-4122)

$elements.Range("AK5").Copy()
$elements.Range("AL5").PasteSpecial(-4122)

$elements.Range("AK6").Copy()
$elements.Range("AL6").PasteSpecial(-4122)
$elements.Range("AL6").Value = "statut"

# Size the new column like the other bestFit text columns
$elements.Columns.Item(38).ColumnWidth = 69.73307291666667
